$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 421.14285
$ws.Range("I9").Value = 495.33334
$ws.Range("K9").Value = 495.33334
$ws.Range("M9").Value = -326.33334

$ws.Range("H100").Value = 3330.4285
$ws.Range("I100").Value = 2256.8
$ws.Range("K100").Value = 2256.8
$ws.Range("M100").Value = -1715.8

$ws.Range("H113").Value = 2538.9285
$ws.Range("I113").Value = 1999.25
$ws.Range("K113").Value = 1999.25
$ws.Range("M113").Value = 1254.75

$ws.Range("H132").Value = 4467.086
$ws.Range("I132").Value = 4635.875
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 13907.625
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -11377.625
$ws.Range("N132").Value = -13060.0001

$ws.Range("H135").Value = 1184.8667
$ws.Range("I135").Value = 1201.6428
$ws.Range("K135").Value = 10814.7852
$ws.Range("M135").Value = -8279.7852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 427
$ws.Range("I5").Value = 491.5
$ws.Range("K5").Value = 491.5
$ws.Range("M5").Value = -379.5

$ws.Range("H16").Value = 2974.5
$ws.Range("I16").Value = 2974.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2974.5
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2687.5

$ws.Range("H21").Value = 13091.3
$ws.Range("I21").Value = 5152.1665
$ws.Range("J21").Value = 25000
$ws.Range("K21").Value = 5152.1665
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = -4778.1665
$ws.Range("N21").Value = -25748

$ws.Range("H45").Value = 3627.5652
$ws.Range("I45").Value = 3138.0557
$ws.Range("K45").Value = 3138.0557
$ws.Range("M45").Value = -2761.0557

$ws.Range("H61").Value = 3090
$ws.Range("I61").Value = 2874.3103
$ws.Range("K61").Value = 2874.3103
$ws.Range("M61").Value = -2662.3103

$ws.Range("H63").Value = 159388.89

$ws.Range("H66").Value = 159388.89

$ws.Range("H74").Value = 3904.4666
$ws.Range("I74").Value = 3904.4666
$ws.Range("K74").Value = 3904.4666
$ws.Range("M74").Value = -3030.4666

$ws.Range("H77").Value = 3904.4666
$ws.Range("I77").Value = 3904.4666
$ws.Range("K77").Value = 19522.333
$ws.Range("M77").Value = -15154.333

$ws.Range("H122").Value = 4395.5713
$ws.Range("I122").Value = 2785.3333
$ws.Range("K122").Value = 8355.999899999999
$ws.Range("M122").Value = -5905.999899999999

$ws.Range("H132").Value = 2685.9048
$ws.Range("I132").Value = 2255.4062
$ws.Range("K132").Value = 6766.2186
$ws.Range("M132").Value = -4236.2186

$ws.Range("H133").Value = 92620.336
$ws.Range("J133").Value = 92620.336
$ws.Range("L133").Value = 92620.336
$ws.Range("N133").Value = -97680.336

$ws.Range("H136").Value = 3090
$ws.Range("I136").Value = 2874.3103
$ws.Range("K136").Value = 8622.930899999999
$ws.Range("M136").Value = -6072.930899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 427
$ws.Range("I4").Value = 491.5
$ws.Range("K4").Value = 491.5
$ws.Range("M4").Value = -376.5

$ws.Range("H10").Value = 2001.3334
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H22").Value = 345
$ws.Range("I22").Value = 345
$ws.Range("K22").Value = 345
$ws.Range("M22").Value = -172

$ws.Range("H126").Value = 106499.5
$ws.Range("J126").Value = 106499.5
$ws.Range("L126").Value = 106499.5
$ws.Range("N126").Value = -116379.5

$ws.Range("H134").Value = 2192.3635
$ws.Range("I134").Value = 2107.8667
$ws.Range("K134").Value = 6323.6001
$ws.Range("M134").Value = -3788.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3189.7144
$ws.Range("I31").Value = 2148.4285
$ws.Range("J31").Value = 4231
$ws.Range("K31").Value = 2148.4285
$ws.Range("L31").Value = 4231
$ws.Range("M31").Value = -1853.4285
$ws.Range("N31").Value = -4821

$ws.Range("H34").Value = 3189.7144
$ws.Range("I34").Value = 2148.4285
$ws.Range("J34").Value = 4231
$ws.Range("K34").Value = 2148.4285
$ws.Range("L34").Value = 4231
$ws.Range("M34").Value = -1946.4285
$ws.Range("N34").Value = -4635

$ws.Range("H58").Value = 2775.7966
$ws.Range("I58").Value = 2530.5952
$ws.Range("J58").Value = 3381.5881
$ws.Range("K58").Value = 2530.5952
$ws.Range("L58").Value = 3381.5881
$ws.Range("M58").Value = -2327.5952
$ws.Range("N58").Value = -3787.5881

$ws.Range("H105").Value = 1867.4615
$ws.Range("I105").Value = 1676.25
$ws.Range("K105").Value = 1676.25
$ws.Range("M105").Value = 70.75

$ws.Range("H132").Value = 1333
$ws.Range("I132").Value = 1333
$ws.Range("K132").Value = 3999
$ws.Range("M132").Value = -1469

$ws.Range("H134").Value = 3014.926
$ws.Range("I134").Value = 2935.55
$ws.Range("J134").Value = 3241.7144
$ws.Range("K134").Value = 8806.650000000001
$ws.Range("L134").Value = 9725.143199999999
$ws.Range("M134").Value = -6271.650000000001
$ws.Range("N134").Value = -14795.1432

$ws.Range("H136").Value = 2775.7966
$ws.Range("I136").Value = 2530.5952
$ws.Range("J136").Value = 3381.5881
$ws.Range("K136").Value = 7591.785600000001
$ws.Range("L136").Value = 10144.7643
$ws.Range("M136").Value = -5041.785600000001
$ws.Range("N136").Value = -15244.7643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4985
$ws.Range("J39").Value = 4985
$ws.Range("L39").Value = 14955
$ws.Range("N39").Value = -15543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2826
$ws.Range("I80").Value = 2152.5
$ws.Range("J80").Value = 3499.5
$ws.Range("K80").Value = 2152.5
$ws.Range("L80").Value = 3499.5
$ws.Range("M80").Value = -1154.5
$ws.Range("N80").Value = -5495.5

$ws.Range("H83").Value = 2826
$ws.Range("I83").Value = 2152.5
$ws.Range("J83").Value = 3499.5
$ws.Range("K83").Value = 10762.5
$ws.Range("L83").Value = 17497.5
$ws.Range("M83").Value = -5770.5
$ws.Range("N83").Value = -27481.5

$ws.Range("H122").Value = 1674.6154
$ws.Range("I122").Value = 1605.8334
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4817.5002
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2367.5002
$ws.Range("N122").Value = -12400

$ws.Range("H132").Value = 4176.375
$ws.Range("I132").Value = 3922.2
$ws.Range("K132").Value = 11766.6
$ws.Range("M132").Value = -9236.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8220
$ws.Range("I7").Value = 8010.706
$ws.Range("K7").Value = 8010.706
$ws.Range("M7").Value = -7898.706

$ws.Range("H100").Value = 5059.2
$ws.Range("I100").Value = 4498.6665
$ws.Range("K100").Value = 4498.6665
$ws.Range("M100").Value = -3957.6665

$ws.Range("H122").Value = 12363.742
$ws.Range("I122").Value = 13257.417
$ws.Range("J122").Value = 9299.714
$ws.Range("K122").Value = 39772.251
$ws.Range("L122").Value = 27899.142
$ws.Range("M122").Value = -37322.251
$ws.Range("N122").Value = -32799.142

$ws.Range("H126").Value = 8220
$ws.Range("I126").Value = 8010.706
$ws.Range("K126").Value = 24032.118
$ws.Range("M126").Value = -21562.118

$ws.Range("H136").Value = 4040.5925
$ws.Range("I136").Value = 3162.7368
$ws.Range("K136").Value = 9488.2104
$ws.Range("M136").Value = -6938.2104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4909.364
$ws.Range("I62").Value = 3700
$ws.Range("K62").Value = 3700
$ws.Range("M62").Value = -3076

$ws.Range("H65").Value = 4909.364
$ws.Range("I65").Value = 3700
$ws.Range("K65").Value = 18500
$ws.Range("M65").Value = -15380

$ws.Range("H100").Value = 1337.4
$ws.Range("I100").Value = 1337.4
$ws.Range("K100").Value = 2674.8
$ws.Range("M100").Value = -2133.8

$ws.Range("H126").Value = 6211.1113
$ws.Range("I126").Value = 6488.2354
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 19464.7062
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -16994.7062
$ws.Range("N126").Value = -9440

$ws.Range("H132").Value = 2571.4211
$ws.Range("I132").Value = 2509.394
$ws.Range("K132").Value = 7528.181999999999
$ws.Range("M132").Value = -4998.181999999999

$ws.Range("H136").Value = 2067.054
$ws.Range("I136").Value = 1424.7693
$ws.Range("J136").Value = 3585.182
$ws.Range("K136").Value = 4274.3079
$ws.Range("L136").Value = 10755.546
$ws.Range("M136").Value = -1724.3079
$ws.Range("N136").Value = -15855.546
